$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 825
$ws1.Range("F9").Value = 1973
$ws1.Range("F11").Value = 353
$ws1.Range("F13").Value = 1612
$ws1.Range("F19").Value = 1444
$ws1.Range("F21").Value = 633
$ws1.Range("F23").Value = 10739
$ws1.Range("F24").Value = 10960

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 1973
$ws4.Range("F13").Value = 353
$ws4.Range("F15").Value = 1612
$ws4.Range("F23").Value = 1444
$ws4.Range("F25").Value = 633
$ws4.Range("F27").Value = 10739
$ws4.Range("F28").Value = 10961
